$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "QuangTrung" user was renamed to "TrungTQ"
$ws.Range("B3").Value = "TrungTQ"

# Drop the stray placeholder avatar text that was sitting in F2
$ws.Range("F2").ClearContents()

# Give the birth-date column a day/month/year display format
$ws.Range("D3").NumberFormat = "dd/mm/yyyy"
